$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.997.97'
$ws.Range('E2').Value = '  +6.46%  '
$ws.Range('D3').Value = '1.883.72'
$ws.Range('E3').Value = '  +5.59%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9997'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '248.59'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.83%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.9993'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4983'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +1.36%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '45.88'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +9.62%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.2855'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +6.56%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.06546'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +4.61%  '
$ws.Range('D11').Value = '1.881.75'
$ws.Range('E11').Value = '  +5.57%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '17.12'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +4.10%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.07206'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +2.48%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.6632'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +5.90%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '85.08'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +6.39%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '4.792'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').Value = '30.003.36'
$ws.Range('E17').Value = '  +6.65%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.10%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.86'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +6.67%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.000007507'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +4.04%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.9992'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = '2.123.40'
$ws.Range('E22').Value = '  +5.73%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.755'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +4.14%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.535'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +5.64%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '9.013'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +3.21%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '145.04'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +2.52%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '134.47'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +23.19%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '16.69'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +5.65%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.958'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +5.23%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.380'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.48%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.172'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.03%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.08612'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.03%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.876'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.66%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.05113'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +4.33%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.131'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +5.48%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.6862'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +5.24%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.705'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.15%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.305'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +12.51%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.755'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +6.39%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.9586'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.34%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.01631'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +5.14%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '6.078'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.19%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '104.11'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +4.30%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.9997'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.08%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.4217'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +5.74%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '7.433'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +3.44%  '
$ws.Range('E48').Value = '  +4.15%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.05634'
$c.Style = 'Normal'
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '32.40'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +5.77%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '8.274'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +3.13%  '
